# Re-apply scheduled market-data refresh to the per-job profit tables (cols H-N)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Cells.Item(38, 8).Value = 943.7778
$ws.Cells.Item(38, 9).Value = 436.75
$ws.Cells.Item(38, 10).Value = 5000
$ws.Cells.Item(38, 11).Value = 1310.25
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 13).Value = -938.25
$ws.Cells.Item(38, 14).Value = -15744

# Row 39
$ws.Cells.Item(39, 8).Value = 20.88889
$ws.Cells.Item(39, 9).Value = 19.75
$ws.Cells.Item(39, 10).Value = 30
$ws.Cells.Item(39, 11).Value = 59.25
$ws.Cells.Item(39, 12).Value = 90
$ws.Cells.Item(39, 13).Value = 236.75
$ws.Cells.Item(39, 14).Value = -682

# Row 64
$ws.Cells.Item(64, 8).Value = 5887.5
$ws.Cells.Item(64, 9).Value = 7500
$ws.Cells.Item(64, 11).Value = 7500
$ws.Cells.Item(64, 13).Value = -7252

# Row 67
$ws.Cells.Item(67, 8).Value = 5887.5
$ws.Cells.Item(67, 9).Value = 7500
$ws.Cells.Item(67, 11).Value = 7500
$ws.Cells.Item(67, 13).Value = -6642

# Row 74
$ws.Cells.Item(74, 8).Value = 1002933.3
$ws.Cells.Item(74, 9).Value = 1002933.3
$ws.Cells.Item(74, 11).Value = 1002933.3
$ws.Cells.Item(74, 13).Value = -1001997.3

# Row 77
$ws.Cells.Item(77, 8).Value = 1002933.3
$ws.Cells.Item(77, 9).Value = 1002933.3
$ws.Cells.Item(77, 11).Value = 5014666.5
$ws.Cells.Item(77, 13).Value = -5009986.5

# Row 80
$ws.Cells.Item(80, 8).Value = 440.58334
$ws.Cells.Item(80, 9).Value = 346.33334
$ws.Cells.Item(80, 10).Value = 534.8333
$ws.Cells.Item(80, 11).Value = 1039.00002
$ws.Cells.Item(80, 12).Value = 1604.4999
$ws.Cells.Item(80, 13).Value = -41.00001999999995
$ws.Cells.Item(80, 14).Value = -3600.4999

# Row 83
$ws.Cells.Item(83, 8).Value = 440.58334
$ws.Cells.Item(83, 9).Value = 346.33334
$ws.Cells.Item(83, 10).Value = 534.8333
$ws.Cells.Item(83, 11).Value = 3117.00006
$ws.Cells.Item(83, 12).Value = 4813.4997
$ws.Cells.Item(83, 13).Value = 1874.99994
$ws.Cells.Item(83, 14).Value = -14797.4997

# Row 132
$ws.Cells.Item(132, 8).Value = 920.2778
$ws.Cells.Item(132, 9).Value = 939.17645
$ws.Cells.Item(132, 11).Value = 2817.52935
$ws.Cells.Item(132, 13).Value = -287.5293500000002

# Row 137
$ws.Cells.Item(137, 8).Value = 1338.325
$ws.Cells.Item(137, 9).Value = 1393.96
$ws.Cells.Item(137, 11).Value = 4181.88
$ws.Cells.Item(137, 13).Value = -1631.88

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5090.5557
$ws.Cells.Item(32, 9).Value = 2493.5833
$ws.Cells.Item(32, 11).Value = 2493.5833
$ws.Cells.Item(32, 13).Value = -2206.5833

# Row 45
$ws.Cells.Item(45, 8).Value = 3377.6667
$ws.Cells.Item(45, 9).Value = 2249.75
$ws.Cells.Item(45, 11).Value = 2249.75
$ws.Cells.Item(45, 13).Value = -1872.75

# Row 74
$ws.Cells.Item(74, 8).Value = 637.94116
$ws.Cells.Item(74, 9).Value = 552.8125
$ws.Cells.Item(74, 11).Value = 552.8125
$ws.Cells.Item(74, 13).Value = 321.1875

# Row 77
$ws.Cells.Item(77, 8).Value = 637.94116
$ws.Cells.Item(77, 9).Value = 552.8125
$ws.Cells.Item(77, 11).Value = 2764.0625
$ws.Cells.Item(77, 13).Value = 1603.9375

# Row 97
$ws.Cells.Item(97, 8).Value = 486.3158
$ws.Cells.Item(97, 9).Value = 466.58823
$ws.Cells.Item(97, 11).Value = 466.58823
$ws.Cells.Item(97, 13).Value = 29.41176999999999

# Row 132
$ws.Cells.Item(132, 8).Value = 2292.9285
$ws.Cells.Item(132, 9).Value = 2008.4166
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 6025.2498
$ws.Cells.Item(132, 12).Value = 12000
$ws.Cells.Item(132, 13).Value = -3495.2498
$ws.Cells.Item(132, 14).Value = -17060

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 2618.4666
$ws.Cells.Item(86, 10).Value = 2050
$ws.Cells.Item(86, 12).Value = 2050
$ws.Cells.Item(86, 14).Value = -4296

# Row 89
$ws.Cells.Item(89, 8).Value = 2618.4666
$ws.Cells.Item(89, 10).Value = 2050
$ws.Cells.Item(89, 12).Value = 10250
$ws.Cells.Item(89, 14).Value = -21482

# Row 134
$ws.Cells.Item(134, 8).Value = 1529.4762
$ws.Cells.Item(134, 9).Value = 1048.3334
$ws.Cells.Item(134, 10).Value = 4416.3335
$ws.Cells.Item(134, 11).Value = 3145.0002
$ws.Cells.Item(134, 12).Value = 13249.0005
$ws.Cells.Item(134, 13).Value = -610.0001999999999
$ws.Cells.Item(134, 14).Value = -18319.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2837.394
$ws.Cells.Item(31, 9).Value = 1381.3334
$ws.Cells.Item(31, 11).Value = 1381.3334
$ws.Cells.Item(31, 13).Value = -1086.3334

# Row 34
$ws.Cells.Item(34, 8).Value = 2837.394
$ws.Cells.Item(34, 9).Value = 1381.3334
$ws.Cells.Item(34, 11).Value = 1381.3334
$ws.Cells.Item(34, 13).Value = -1179.3334

# Row 132
$ws.Cells.Item(132, 8).Value = 2231.9333
$ws.Cells.Item(132, 9).Value = 1735.125
$ws.Cells.Item(132, 11).Value = 5205.375
$ws.Cells.Item(132, 13).Value = -2675.375

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Cells.Item(44, 8).Value = 898.6667
$ws.Cells.Item(44, 9).Value = 319.14285
$ws.Cells.Item(44, 11).Value = 957.4285500000001
$ws.Cells.Item(44, 13).Value = -559.4285500000001

# Row 112
$ws.Cells.Item(112, 8).Value = 3166.6667

$ws = $wb.Worksheets.Item("GSM")
# Row 114
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()

# Row 123
$ws.Cells.Item(123, 8).Value = 31490.3
$ws.Cells.Item(123, 10).Value = 31111.445
$ws.Cells.Item(123, 12).Value = 31111.445
$ws.Cells.Item(123, 14).Value = -36011.445

# Row 132
$ws.Cells.Item(132, 8).Value = 2203.0667
$ws.Cells.Item(132, 10).Value = 3249.75
$ws.Cells.Item(132, 12).Value = 9749.25
$ws.Cells.Item(132, 14).Value = -14809.25

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Cells.Item(9, 8).Value = 696.6667
$ws.Cells.Item(9, 9).Value = 45
$ws.Cells.Item(9, 11).Value = 45
$ws.Cells.Item(9, 13).Value = 179

# Row 19
$ws.Cells.Item(19, 8).Value = 10499.667
$ws.Cells.Item(19, 9).Value = 1500
$ws.Cells.Item(19, 10).Value = 14999.5
$ws.Cells.Item(19, 11).Value = 1500
$ws.Cells.Item(19, 12).Value = 14999.5
$ws.Cells.Item(19, 13).Value = -1330
$ws.Cells.Item(19, 14).Value = -15339.5

# Row 22
$ws.Cells.Item(22, 8).Value = 5150
$ws.Cells.Item(22, 10).Value = 3767
$ws.Cells.Item(22, 12).Value = 3767
$ws.Cells.Item(22, 14).Value = -4357

# Row 27
$ws.Cells.Item(27, 8).Value = 5150
$ws.Cells.Item(27, 10).Value = 3767
$ws.Cells.Item(27, 12).Value = 3767
$ws.Cells.Item(27, 14).Value = -3981

# Row 132
$ws.Cells.Item(132, 8).Value = 4197.5386
$ws.Cells.Item(132, 9).Value = 4006.3
$ws.Cells.Item(132, 10).Value = 4835
$ws.Cells.Item(132, 11).Value = 12018.9
$ws.Cells.Item(132, 12).Value = 14505
$ws.Cells.Item(132, 13).Value = -9488.900000000001
$ws.Cells.Item(132, 14).Value = -19565

# Row 133
$ws.Cells.Item(133, 8).Value = 60000
$ws.Cells.Item(133, 10).Value = 60000
$ws.Cells.Item(133, 12).Value = 60000
$ws.Cells.Item(133, 14).Value = -65060

# Row 136
$ws.Cells.Item(136, 8).Value = 5077.04
$ws.Cells.Item(136, 9).Value = 4782.7896
$ws.Cells.Item(136, 10).Value = 6008.8335
$ws.Cells.Item(136, 11).Value = 14348.3688
$ws.Cells.Item(136, 12).Value = 18026.5005
$ws.Cells.Item(136, 13).Value = -11798.3688
$ws.Cells.Item(136, 14).Value = -23126.5005

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 1216.6666
$ws.Cells.Item(81, 9).Value = 1370
$ws.Cells.Item(81, 10).Value = 450
$ws.Cells.Item(81, 11).Value = 2740
$ws.Cells.Item(81, 12).Value = 900
$ws.Cells.Item(81, 13).Value = -1679
$ws.Cells.Item(81, 14).Value = -3022

# Row 84
$ws.Cells.Item(84, 8).Value = 1216.6666
$ws.Cells.Item(84, 9).Value = 1370
$ws.Cells.Item(84, 10).Value = 450
$ws.Cells.Item(84, 11).Value = 13700
$ws.Cells.Item(84, 12).Value = 4500
$ws.Cells.Item(84, 13).Value = -8396
$ws.Cells.Item(84, 14).Value = -15108

# Row 107
$ws.Cells.Item(107, 8).Value = 438.66666
$ws.Cells.Item(107, 9).Value = 443.625
$ws.Cells.Item(107, 11).Value = 1330.875
$ws.Cells.Item(107, 13).Value = 589.125

# Row 108
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()

# Row 132
$ws.Cells.Item(132, 8).Value = 2169.125
$ws.Cells.Item(132, 9).Value = 1665.8823
$ws.Cells.Item(132, 11).Value = 4997.6469
$ws.Cells.Item(132, 13).Value = -2467.6469

# Row 136
$ws.Cells.Item(136, 8).Value = 1478.95
$ws.Cells.Item(136, 9).Value = 848.75
$ws.Cells.Item(136, 11).Value = 2546.25
$ws.Cells.Item(136, 13).Value = 3.75
